# Edit: insert two new price records (rows 814-815) into the "Papa" sheet,
# pushing the existing rows 814-880 down to 816-882.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 814; everything from old row 814 onward
# shifts down by two rows (old 814 -> new 816, ..., old 880 -> new 882).
$ws.Range("A814:A815").EntireRow.Insert()

# Populate new row 814
$ws.Cells.Item(814, 1).Value = 9
$ws.Cells.Item(814, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(814, 3).Value = "Metropolitana"
$ws.Cells.Item(814, 4).Value = 44578
$ws.Cells.Item(814, 5).Value = 13
$ws.Cells.Item(814, 6).Value = 100114001
$ws.Cells.Item(814, 7).Value = "Papa"
$ws.Cells.Item(814, 8).Value = "Asterix"
$ws.Cells.Item(814, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(814, 10).Value = 340
$ws.Cells.Item(814, 11).Value = 11000
$ws.Cells.Item(814, 12).Value = 12000
$ws.Cells.Item(814, 13).Value = 11500
$ws.Cells.Item(814, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(814, 15).Value = "Región del Maule"
$ws.Cells.Item(814, 16).Value = 460
$ws.Cells.Item(814, 17).Value = 25
$ws.Cells.Item(814, 18).Value = "Hortaliza"

# Populate new row 815
$ws.Cells.Item(815, 1).Value = 9
$ws.Cells.Item(815, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(815, 3).Value = "Metropolitana"
$ws.Cells.Item(815, 4).Value = 44578
$ws.Cells.Item(815, 5).Value = 13
$ws.Cells.Item(815, 6).Value = 100114001
$ws.Cells.Item(815, 7).Value = "Papa"
$ws.Cells.Item(815, 8).Value = "Asterix"
$ws.Cells.Item(815, 9).Value = "1a (cosecha)"
$ws.Cells.Item(815, 10).Value = 160
$ws.Cells.Item(815, 11).Value = 9000
$ws.Cells.Item(815, 12).Value = 10000
$ws.Cells.Item(815, 13).Value = 9500
$ws.Cells.Item(815, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(815, 15).Value = "Región del Maule"
$ws.Cells.Item(815, 16).Value = 380
$ws.Cells.Item(815, 17).Value = 25
$ws.Cells.Item(815, 18).Value = "Hortaliza"
